$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pre-existing bordered formatting on A1:B1 before writing the
# appended header/data rows so the new cells pick up the default style.
$ws.Range("A1:B1").ClearFormats()

# Header row
$ws.Range("A1").Value = "MIGRATION DATE"
$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"

# Data row - the migration date is stored as literal text, not a date
# serial, so force text entry with a leading apostrophe and then drop the
# resulting "quote prefix" style so the cell keeps the default format.
$ws.Range("A2").Value = "'2025-10-16"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "YYY"
$ws.Range("C2").Value = "123ABX007"
$ws.Range("D2").Value = "Karapakkam"
